# Regenerate save_data column G ("K") with freshly simulated strike-count
# values (commit: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals").
#
# Column G holds the new simulated "K" values for each trade row (rows 2-72,
# column 7). The values below are the freshly (re)computed s_vals for this
# run, keyed by worksheet row number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2  = 1
    3  = 0
    4  = 1
    5  = 3
    6  = 1
    7  = 2
    8  = 0
    9  = 0
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 2
    19 = 1
    20 = 1
    21 = 0
    22 = 2
    23 = 2
    24 = 2
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 3
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 1
    39 = 2
    40 = 0
    41 = 0
    42 = 2
    43 = 1
    44 = 0
    45 = 1
    46 = 3
    47 = 3
    48 = 1
    49 = 1
    50 = 2
    51 = 3
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 2
    57 = 1
    58 = 0
    59 = 3
    60 = 0
    61 = 2
    62 = 0
    63 = 1
    64 = 3
    65 = 1
    66 = 0
    67 = 0
    68 = 1
    69 = 3
    70 = 1
    72 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
